# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" on all sheets
# - Narrow the "Status" column(s) from ~17.22 chars to ~13.41 chars

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status value ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns("E:F").ColumnWidth = 13.41

# --- zh-cn sheet: column C ("Status") ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns("C:C").ColumnWidth = 13.41

# --- de-de sheet: column C ("Status") ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns("C:C").ColumnWidth = 13.41
